# Advanced Reactor Materials / grading.xlsx — "updated teaching and grades"
#
# Adds Quiz #4 scores + bonus-adjusted Quiz #4 (J/K), a Final Project column
# (N) with a Penalty column (O), shifts the old "Bonus" header + running
# total from column O/P to P/Q, adds a grade-scale note, and recomputes the
# dependent weighting/running-total formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row 7: old "Bonus" (O7) becomes "Penalty"; "Bonus" moves to P7
# ---------------------------------------------------------------------
$ws.Range("O7").Value2 = "Penalty"
$ws.Range("P7").Value2 = "Bonus"

# ---------------------------------------------------------------------
# Quiz #4 (J) + Quiz #4 * (K, bonus-adjusted) for each student
# ---------------------------------------------------------------------
$ws.Range("J8").Value2 = 93
$ws.Range("K8").Formula = "=J8+5"

$ws.Range("J9").Value2 = 90
$ws.Range("K9").Formula = "=J9+5"

$ws.Range("J10").Value2 = 96
$ws.Range("K10").Formula = "=J10+5"

# ---------------------------------------------------------------------
# Final Project (N, red text) + Penalty (O) for each student
# ---------------------------------------------------------------------
$ws.Range("N8:N10").Font.Color = 255
$ws.Range("N8").Value2 = 95
$ws.Range("N9").Value2 = 95
$ws.Range("N10").Value2 = 95

$ws.Range("O8").Value2 = 5
$ws.Range("O9").Value2 = 10
$ws.Range("O10").Value2 = 5

# ---------------------------------------------------------------------
# Grade scale note (Times New Roman) next to rows 9 & 10
# ---------------------------------------------------------------------
$ws.Range("Q9:Q10").Font.Name = "Times New Roman"
$ws.Range("Q9").Value2 = "A+ 98-100; A 93-97; A- 90-92"
$ws.Range("Q10").Value2 = "B+ 87-89; B 83-87; B- 80-82"

# ---------------------------------------------------------------------
# Row 12 averages now cover the new J/K/N columns as well
# ---------------------------------------------------------------------
$ws.Range("D12:I12").Formula = "=AVERAGE(D8:D10)"
$ws.Range("J12:N12").Formula = "=AVERAGE(J8:J10)"
$ws.Range("K12").Formula = "=AVERAGE(K8:K10)"

# ---------------------------------------------------------------------
# Row 13 weights: new Penalty weight (O13); old running-total cell (P13)
# is retired now that the running total lives under column Q
# ---------------------------------------------------------------------
$ws.Range("O13").Value2 = 20
$ws.Range("P13").Clear()

# ---------------------------------------------------------------------
# Row 15 "Running Total" label moves from P15 to Q15
# ---------------------------------------------------------------------
$ws.Range("P15").Clear()
$ws.Range("Q15").Value2 = "Running Total"

# ---------------------------------------------------------------------
# Row 16 max-points: new K16/O16 entries; running total now in Q16
# ---------------------------------------------------------------------
$ws.Range("K16").Value2 = 12.5
$ws.Range("O16").Value2 = 20
$ws.Range("P16").Clear()
$ws.Range("Q16").Formula = "=SUM(D16:O16)"

# ---------------------------------------------------------------------
# Rows 17-19: per-student earned points, now including the Penalty-adjusted
# Final Project (O) column; percentage total moves from P to Q
# ---------------------------------------------------------------------
$ws.Range("O17").Formula = "=(N8-O8)*O`$13/100"
$ws.Range("P17").Clear()
$ws.Range("Q17").NumberFormat = "0.00"
$ws.Range("Q17").Formula = "=SUM(D17:O17)/`$Q`$16"

$ws.Range("O18").Formula = "=(N9-O9)*O`$13/100"
$ws.Range("P18").Clear()
$ws.Range("Q18").NumberFormat = "0.00"
$ws.Range("Q18").Formula = "=SUM(D18:O18)/`$Q`$16"

$ws.Range("O19").Formula = "=(N10-O10)*O`$13/100"
$ws.Range("P19").Clear()
$ws.Range("Q19").NumberFormat = "0.00"
$ws.Range("Q19").Formula = "=SUM(D19:O19)/`$Q`$16"

# ---------------------------------------------------------------------
# Selection cursor, matching the saved workbook view
# ---------------------------------------------------------------------
$ws.Range("R12").Select()
